# Add support for reading empty cell and row.
# This adds a new row of sample data (row 4), skips row 5 on purpose
# (to exercise "empty row" handling), and adds two more rows (6 and 7)
# separated from row 4 by the empty row 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: plain UTF-8 / multi-byte text in column A and C (column B is
# intentionally left empty to exercise "empty cell" handling), plus
# bold / italic / underline styled cells in F, G, H.
$ws.Range("A4").Value = "utf"
$ws.Range("C4").Value = "日本語"

$ws.Range("F4").Value = "bold"
$ws.Range("F4").Font.Bold = $true

$ws.Range("G4").Value = "italic"
$ws.Range("G4").Font.Italic = $true

$ws.Range("H4").Value = "underline"
$ws.Range("H4").Font.Underline = $true

# Row 5 is intentionally left blank (empty row) between rows 4 and 6.

# Rows 6 and 7: simple single-cell rows after the empty row.
$ws.Range("A6").Value = "Jump Row"
$ws.Range("A7").Value = "Another Row"

# Move the active selection past the new data, like the source workbook.
[void]$ws.Range("A8").Select()

# Configure the page setup (paper size 9 = A4, portrait orientation).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
